$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 369
$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # Column C: "Förändrad" date bumps from 45184 to 45186
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value2 = 45186
    }

    # Hyperlink formulas (S, T, U, V, W, X, Y) gain a friendly-name second argument
    # equal to the row's "Beteckning" (column A) value.
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ($beteckning) {
        foreach ($col in $linkCols) {
            $cell = $ws.Range($col + $r)
            $f = $cell.Formula
            if ($f -like '=HYPERLINK(*' -and $f -notlike '*,*') {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
